# Scheduled runner update: refresh market-price derived columns (H-N)
# across the leve-profit sheets, per latest Universalis price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2905707.8
$ws.Range("J17").Value = 3005687.2
$ws.Range("L17").Value = 9017061.600000001
$ws.Range("N17").Value = -9017397.600000001
$ws.Range("H28").Value = 548.55554
$ws.Range("I28").Value = 476.27274
$ws.Range("K28").Value = 476.27274
$ws.Range("M28").Value = 8.727260000000001
$ws.Range("H33").Value = 1425112.9
$ws.Range("I33").Value = 1851664.4
$ws.Range("K33").Value = 1851664.4
$ws.Range("M33").Value = -1851435.4
$ws.Range("H39").Value = 545.7143
$ws.Range("I39").Value = 424
$ws.Range("K39").Value = 1272
$ws.Range("M39").Value = -976
$ws.Range("H100").Value = 47145.41
$ws.Range("I100").Value = 63468.812
$ws.Range("J100").Value = 3616.3333
$ws.Range("K100").Value = 63468.812
$ws.Range("L100").Value = 3616.3333
$ws.Range("M100").Value = -62927.812
$ws.Range("N100").Value = -4698.3333
$ws.Range("H132").Value = 2418.0278
$ws.Range("I132").Value = 2265.9707
$ws.Range("K132").Value = 6797.9121
$ws.Range("M132").Value = -4267.9121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3644.6528
$ws.Range("I32").Value = 3644.6528
$ws.Range("K32").Value = 3644.6528
$ws.Range("M32").Value = -3357.6528
$ws.Range("H43").Value = 5397.5
$ws.Range("H45").Value = 6355.9165
$ws.Range("I45").Value = 7574.8887
$ws.Range("K45").Value = 7574.8887
$ws.Range("M45").Value = -7197.8887
$ws.Range("H60").Value = 9999
$ws.Range("I60").Value = 9999
$ws.Range("K60").Value = 9999
$ws.Range("M60").Value = -9266
$ws.Range("H97").Value = 1013.94446
$ws.Range("I97").Value = 902.4375
$ws.Range("K97").Value = 902.4375
$ws.Range("M97").Value = -406.4375
$ws.Range("H122").Value = 3227.75
$ws.Range("I122").Value = 2974.5715
$ws.Range("K122").Value = 8923.7145
$ws.Range("M122").Value = -6473.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 670622.9399999999
$ws.Range("I86").Value = 2001880.8
$ws.Range("K86").Value = 2001880.8
$ws.Range("M86").Value = -2000757.8
$ws.Range("H89").Value = 670622.9399999999
$ws.Range("I89").Value = 2001880.8
$ws.Range("K89").Value = 10009404
$ws.Range("M89").Value = -10003788
$ws.Range("H99").Value = 3338.018
$ws.Range("I99").Value = 3349.1914
$ws.Range("J99").Value = 3272.375
$ws.Range("K99").Value = 3349.1914
$ws.Range("L99").Value = 3272.375
$ws.Range("M99").Value = -1851.1914
$ws.Range("N99").Value = -6268.375
$ws.Range("H105").Value = 3715.8215
$ws.Range("I105").Value = 3301.5908
$ws.Range("K105").Value = 3301.5908
$ws.Range("M105").Value = -1554.5908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 9466
$ws.Range("J8").Value = 9466
$ws.Range("L8").Value = 9466
$ws.Range("N8").Value = -9746
$ws.Range("H31").Value = 4143.5713
$ws.Range("I31").Value = 2501.25
$ws.Range("J31").Value = 6333.3335
$ws.Range("K31").Value = 2501.25
$ws.Range("L31").Value = 6333.3335
$ws.Range("M31").Value = -2206.25
$ws.Range("N31").Value = -6923.3335
$ws.Range("H34").Value = 4143.5713
$ws.Range("I34").Value = 2501.25
$ws.Range("J34").Value = 6333.3335
$ws.Range("K34").Value = 2501.25
$ws.Range("L34").Value = 6333.3335
$ws.Range("M34").Value = -2299.25
$ws.Range("N34").Value = -6737.3335
$ws.Range("H58").Value = 2857.5693
$ws.Range("I58").Value = 1946.1
$ws.Range("K58").Value = 1946.1
$ws.Range("M58").Value = -1743.1
$ws.Range("H132").Value = 23190.639
$ws.Range("I132").Value = 12156.846
$ws.Range("K132").Value = 36470.538
$ws.Range("M132").Value = -33940.538
$ws.Range("H136").Value = 2857.5693
$ws.Range("I136").Value = 1946.1
$ws.Range("K136").Value = 5838.299999999999
$ws.Range("M136").Value = -3288.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2409
$ws.Range("I14").Value = 2409
$ws.Range("K14").Value = 7227
$ws.Range("M14").Value = -7054
$ws.Range("H29").Value = 2514
$ws.Range("I29").Value = 28
$ws.Range("J29").Value = 5000
$ws.Range("K29").Value = 84
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = 193
$ws.Range("N29").Value = -15554
$ws.Range("H98").Value = 1373.75
$ws.Range("I98").Value = 1376.25
$ws.Range("J98").Value = 1372.0834
$ws.Range("K98").Value = 4128.75
$ws.Range("L98").Value = 4116.2502
$ws.Range("M98").Value = -2630.75
$ws.Range("N98").Value = -7112.2502
$ws.Range("H121").Value = 3054.4482
$ws.Range("J121").Value = 3316.8462
$ws.Range("L121").Value = 9950.5386
$ws.Range("N121").Value = -12570.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 4500
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H80").Value = 2799.6667
$ws.Range("I80").Value = 2799.6667
$ws.Range("K80").Value = 2799.6667
$ws.Range("M80").Value = -1801.6667
$ws.Range("H83").Value = 2799.6667
$ws.Range("I83").Value = 2799.6667
$ws.Range("K83").Value = 13998.3335
$ws.Range("M83").Value = -9006.333500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8480.762000000001
$ws.Range("I16").Value = 1070.1875
$ws.Range("K16").Value = 1070.1875
$ws.Range("M16").Value = -900.1875
$ws.Range("H55").Value = 284.58823
$ws.Range("I55").Value = 220.22223
$ws.Range("J55").Value = 357
$ws.Range("K55").Value = 220.22223
$ws.Range("L55").Value = 357
$ws.Range("M55").Value = -47.22223
$ws.Range("N55").Value = -703
$ws.Range("H61").Value = 1258.5
$ws.Range("I61").Value = 1241
$ws.Range("J61").Value = 1299.3334
$ws.Range("K61").Value = 1241
$ws.Range("L61").Value = 1299.3334
$ws.Range("M61").Value = -1039
$ws.Range("N61").Value = -1703.3334
$ws.Range("H74").Value = 30745.273
$ws.Range("I74").Value = 24599.857
$ws.Range("J74").Value = 41499.75
$ws.Range("K74").Value = 24599.857
$ws.Range("L74").Value = 41499.75
$ws.Range("M74").Value = -23601.857
$ws.Range("N74").Value = -43495.75
$ws.Range("H77").Value = 30745.273
$ws.Range("I77").Value = 24599.857
$ws.Range("J77").Value = 41499.75
$ws.Range("K77").Value = 73799.571
$ws.Range("L77").Value = 124499.25
$ws.Range("M77").Value = -68807.571
$ws.Range("N77").Value = -134483.25
$ws.Range("H113").Value = 1258.5
$ws.Range("I113").Value = 1241
$ws.Range("J113").Value = 1299.3334
$ws.Range("K113").Value = 1241
$ws.Range("L113").Value = 1299.3334
$ws.Range("M113").Value = 929
$ws.Range("N113").Value = -5639.3334
$ws.Range("H134").Value = 97214.39999999999
$ws.Range("J134").Value = 97214.39999999999
$ws.Range("L134").Value = 97214.39999999999
$ws.Range("N134").Value = -107354.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 31770.889
$ws.Range("I61").Value = 31770.889
$ws.Range("K61").Value = 31770.889
$ws.Range("M61").Value = -31478.889
$ws.Range("H81").Value = 1621.238
$ws.Range("J81").Value = 2888.889
$ws.Range("L81").Value = 5777.778
$ws.Range("N81").Value = -7899.778
$ws.Range("H84").Value = 1621.238
$ws.Range("J84").Value = 2888.889
$ws.Range("L84").Value = 28888.89
$ws.Range("N84").Value = -39496.89
$ws.Range("H110").Value = 94994
$ws.Range("J110").Value = 94994
$ws.Range("L110").Value = 94994
$ws.Range("N110").Value = -103174
$ws.Range("H126").Value = 8741.200000000001
$ws.Range("I126").Value = 5834.577
$ws.Range("J126").Value = 17138.111
$ws.Range("K126").Value = 17503.731
$ws.Range("L126").Value = 51414.333
$ws.Range("M126").Value = -15033.731
$ws.Range("N126").Value = -56354.333
